# The source edit replaces the empty paragraph that sits right after the
# "Objetivo" heading (the very last paragraph of the document body) with a
# paragraph holding two runs of text.
$d = $word.ActiveDocument

# Locate the "Objetivo" heading paragraph by its text so the script does
# not depend on a brittle, hard-coded paragraph index.
$objetivoIndex = -1
$count = $d.Content.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $candidateText = $d.Content.Paragraphs.Item($i).Range.Text.Trim()
    if ($candidateText -eq "Objetivo") {
        $objetivoIndex = $i
    }
}

if ($objetivoIndex -eq -1) {
    throw "Could not find the 'Objetivo' paragraph"
}

$targetIndex = $objetivoIndex + 1
$targetRange = $d.Content.Paragraphs.Item($targetIndex).Range

# Type the first sentence into the existing (empty) paragraph; it becomes
# that paragraph's first run.
$targetRange.InsertAfter("Producir código basados en los principios de la programación orientada a objetos, que éste sea descriptivo ")

# Split off a new paragraph and type the second sentence into it, so it is
# created as a separate run rather than being coalesced into the first
# sentence's run.
$targetRange.InsertParagraphAfter()
$secondRange = $d.Content.Paragraphs.Item($targetIndex + 1).Range
$secondRange.InsertAfter("sin necesidad de comentar cada linea, haciendo un desarrollo basándonos en las buenas prácticas con la ayuda de herramientas como GIT, GITHUB, MAVEN, entre otras.")

# Join the two paragraphs back together by deleting the paragraph mark
# between them, leaving a single paragraph that contains both runs - the
# exact shape the diff calls for.
$firstRange = $d.Content.Paragraphs.Item($targetIndex).Range
$paraMark = $d.Range($firstRange.End - 1, $firstRange.End)
$paraMark.Delete()

Write-Output "Objetivo paragraph updated"
